{"js": "// Apply the 5 placeholder-renaming edits described in the diff.\n// Each target paragraph holds exactly one run/one <w:t>, so we can\n// safely replace the whole paragraph's text in one shot.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst replacements = [\n  [\"\u03a3\u03c4\u03b7\u03bd {{ place }} \u03c3\u03ae\u03bc\u03b5\u03c1\u03b1 \u03c4\u03b7\u03bd {{ date_num_apologia}} \u03c4\u03bf\u03c5 \u03bc\u03ae\u03bd\u03b1 {{ month_apologia }} \u03c4\u03bf\u03c5 \u03ad\u03c4\u03bf\u03c5\u03c2 {{year_apologia}} \u03b7\u03bc\u03ad\u03c1\u03b1 \u03b5\u03b2\u03b4\u03bf\u03bc\u03ac\u03b4\u03b1\u03c2 {{ day_apologia}} \u03ba\u03b1\u03b9 \u03ce\u03c1\u03b1 {{start_hour_apologia}} \u03b5\u03bd\u03ce\u03c0\u03b9\u03bf\u03bd \u03b5\u03bc\u03bf\u03cd \u03c4\u03bf\u03c5 {{ first_officer }} \u03c4\u03bf\u03c5 {{ policeStation }} \u0398\u03b5\u03c3\u03c3\u03b1\u03bb\u03bf\u03bd\u03af\u03ba\u03b7\u03c2, \u03c0\u03b1\u03c1\u03b9\u03c3\u03c4\u03b1\u03bc\u03ad\u03bd\u03bf\u03c5  \u03ba\u03b1\u03b9 \u03c4\u03bf\u03c5  {{ sec_officer }} \u03c4\u03b7\u03c2 \u03b9\u03b4\u03af\u03b1\u03c2 \u03c5\u03c0\u03b7\u03c1\u03b5\u03c3\u03af\u03b1\u03c2, \u03c0\u03bf\u03c5 \u03c0\u03c1\u03bf\u03c3\u03bb\u03ae\u03c6\u03b8\u03b7\u03ba\u03b5 \u03c9\u03c2 \u0392' \u0391\u03bd\u03b1\u03ba\u03c1\u03b9\u03c4\u03b9\u03ba\u03cc\u03c2 \u03a5\u03c0\u03ac\u03bb\u03bb\u03b7\u03bb\u03bf\u03c2, \u03b5\u03be\u03b5\u03c4\u03ac\u03b6\u03b5\u03c4\u03b1\u03b9 \u03bf \u03ba\u03b1\u03c4\u03c9\u03c4\u03ad\u03c1\u03c9 \u03c3\u03b7\u03bc\u03b5\u03b9\u03bf\u03cd\u03bc\u03b5\u03bd\u03bf\u03c2 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03bf\u03cd\u03bc\u03b5\u03bd\u03bf\u03c2 ,\", \"\u03a3\u03c4\u03b7\u03bd {{ place1 }} \u03c3\u03ae\u03bc\u03b5\u03c1\u03b1 \u03c4\u03b7\u03bd {{ date_num_apologia1}} \u03c4\u03bf\u03c5 \u03bc\u03ae\u03bd\u03b1 {{ month_apologia1 }} \u03c4\u03bf\u03c5 \u03ad\u03c4\u03bf\u03c5\u03c2 {{year_apologia1}} \u03b7\u03bc\u03ad\u03c1\u03b1 \u03b5\u03b2\u03b4\u03bf\u03bc\u03ac\u03b4\u03b1\u03c2 {{day_apologia1}} \u03ba\u03b1\u03b9 \u03ce\u03c1\u03b1 {{start_hour_apologia1}} \u03b5\u03bd\u03ce\u03c0\u03b9\u03bf\u03bd \u03b5\u03bc\u03bf\u03cd \u03c4\u03bf\u03c5 {{ first_officer }} \u03c4\u03bf\u03c5 {{ policeStation }} \u0398\u03b5\u03c3\u03c3\u03b1\u03bb\u03bf\u03bd\u03af\u03ba\u03b7\u03c2, \u03c0\u03b1\u03c1\u03b9\u03c3\u03c4\u03b1\u03bc\u03ad\u03bd\u03bf\u03c5  \u03ba\u03b1\u03b9 \u03c4\u03bf\u03c5  {{ sec_officer }} \u03c4\u03b7\u03c2 \u03b9\u03b4\u03af\u03b1\u03c2 \u03c5\u03c0\u03b7\u03c1\u03b5\u03c3\u03af\u03b1\u03c2, \u03c0\u03bf\u03c5 \u03c0\u03c1\u03bf\u03c3\u03bb\u03ae\u03c6\u03b8\u03b7\u03ba\u03b5 \u03c9\u03c2 \u0392' \u0391\u03bd\u03b1\u03ba\u03c1\u03b9\u03c4\u03b9\u03ba\u03cc\u03c2 \u03a5\u03c0\u03ac\u03bb\u03bb\u03b7\u03bb\u03bf\u03c2, \u03b5\u03be\u03b5\u03c4\u03ac\u03b6\u03b5\u03c4\u03b1\u03b9 \u03bf \u03ba\u03b1\u03c4\u03c9\u03c4\u03ad\u03c1\u03c9 \u03c3\u03b7\u03bc\u03b5\u03b9\u03bf\u03cd\u03bc\u03b5\u03bd\u03bf\u03c2 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03bf\u03cd\u03bc\u03b5\u03bd\u03bf\u03c2 ,\"],\n  [\"\u0391\u03a0\u039f\u039a\u03a1\u0399\u03a3\u0397 : {{surnamePerperator}}  {{namePerperator}} \u03c4\u03bf\u03c5 {{fathernamePerperator}} \u03ba\u03b1\u03b9 \u03c4\u03b7\u03c2 {{ mothernamePerperator}} \u03b3\u03b5\u03bd. {{dateOfBirthPerperator }} \u03c3\u03c4\u03b7 {{ placeOfBirthPerperator }} \u03ba\u03b1\u03c4.{{ addressPerperator }},\u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03c2 \u03c4\u03b7\u03bb\u03b5\u03c6\u03ce\u03bd\u03bf\u03c5 {{ telPreperator }}, \u03b7\u03bb\u03b5\u03ba\u03c4\u03c1\u03bf\u03bd\u03b9\u03ba\u03cc \u03c4\u03b1\u03c7\u03c5\u03b4\u03c1\u03bf\u03bc\u03b5\u03af\u03bf\u03c5{{ emailPreperator }}, \u03ba\u03ac\u03c4\u03bf\u03c7\u03bf\u03c2 \u03c4\u03bf\u03c5 \u03c5\u03c0 \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03bd {{ DATperperator }} \u03c0\u03bf\u03c5 \u03b5\u03ba\u03b4\u03cc\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd {{ issuedPerperator }} \u03b1\u03c0\u03cc {{place_issuedPerperator}}\u0391.\u03a6.\u039c : {{ afmPreperator }}, \u0394.\u039f.\u03a5 : {{ doyPrep }}, \u03b1\u03c0\u03cc \u03c4\u03bf\u03bd {{officer_arrest }} \u03c5\u03c0\u03b7\u03c1\u03b5\u03c4\u03ce\u03bd \u03c3\u03c4\u03bf  {{  policeStation }}\", \"\u0391\u03a0\u039f\u039a\u03a1\u0399\u03a3\u0397 : {{surname}} {{name}} \u03c4\u03bf\u03c5 {{ fathername}} \u03ba\u03b1\u03b9 \u03c4\u03b7\u03c2 {{ mothername}} \u03b3\u03b5\u03bd. {{dateOfBirth  }} \u03c3\u03c4\u03b7 {{ placeOfBirth }} \u03ba\u03b1\u03c4. {{ address }}, \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03c2 \u03c4\u03b7\u03bb\u03b5\u03c6\u03ce\u03bd\u03bf\u03c5 {{ tel }}, \u03b7\u03bb\u03b5\u03ba\u03c4\u03c1\u03bf\u03bd\u03b9\u03ba\u03cc \u03c4\u03b1\u03c7\u03c5\u03b4\u03c1\u03bf\u03bc\u03b5\u03af\u03bf {{ email }}, \u03ba\u03ac\u03c4\u03bf\u03c7\u03bf\u03c2 \u03c4\u03bf\u03c5 \u03c5\u03c0 \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03bd {{ DAT }} \u03c0\u03bf\u03c5 \u03b5\u03ba\u03b4\u03cc\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd  {{ issued }} \u03b1\u03c0\u03cc {{place_issued}} \u0391.\u03a6.\u039c : {{ afm }}, \u0394.\u039f.\u03a5 : {{ doy }}, \u03b1\u03c0\u03cc \u03c4\u03bf\u03bd {{officer_arrest }} \u03c5\u03c0\u03b7\u03c1\u03b5\u03c4\u03ce\u03bd \u03c3\u03c4\u03bf  {{  policeStation }}\"],\n  [\"\u0395\u03bd\u03c4\u03b1\u03cd\u03b8\u03b1 \u03b3\u03bd\u03c9\u03c1\u03af\u03c3\u03b1\u03bc\u03b5 \u03c3\u03c4\u03bf\u03bd \u03b5\u03be\u03b5\u03c4\u03b1\u03b6\u03cc\u03bc\u03b5\u03bd\u03bf \u03cc\u03c4\u03b9 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03b5\u03af\u03c4\u03b1\u03b9 \u03b3\u03b9\u03b1 \u03c0\u03b1\u03c1\u03ac\u03b2\u03b1\u03c3\u03b7 \u03c4\u03bf\u03c5 \u03ac\u03c1\u03b8\u03c1\u03bf\u03c5 {{offences }} \u03c4\u03bf\u03c5 \u03a0.\u039a.\", \"\u0395\u03bd\u03c4\u03b1\u03cd\u03b8\u03b1 \u03b3\u03bd\u03c9\u03c1\u03af\u03c3\u03b1\u03bc\u03b5 \u03c3\u03c4\u03bf\u03bd \u03b5\u03be\u03b5\u03c4\u03b1\u03b6\u03cc\u03bc\u03b5\u03bd\u03bf \u03cc\u03c4\u03b9 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03b5\u03af\u03c4\u03b1\u03b9 \u03b3\u03b9\u03b1 \u03c0\u03b1\u03c1\u03ac\u03b2\u03b1\u03c3\u03b7 \u03c4\u03bf\u03c5 \u03ac\u03c1\u03b8\u03c1\u03bf\u03c5 {{offences_text1}} \u03c4\u03bf\u03c5 \u03a0.\u039a.\"],\n  [\"\u0395\u03b9\u03b4\u03b9\u03ba\u03cc\u03c4\u03b5\u03c1\u03b1 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03b5\u03af\u03c4\u03b1\u03b9 \u03cc\u03c4\u03b9 \u03c3\u03c4\u03b9\u03c2 {{dateOfCrime}}\u03ba\u03b1\u03b9 \u03c0\u03b5\u03c1\u03af \u03ce\u03c1\u03b1 {{ hourOfCrime}} \u03c3\u03c4\u03bf {{placeOfCrime}}  \u03c0\u03c1\u03bf\u03ad\u03b2\u03b7\u03c2 \u03c3\u03c4\u03b7\u03bd \u03c0\u03b1\u03c1\u03ac\u03b2\u03b1\u03c3\u03b7 \u03c4\u03bf\u03c5/\u03c4\u03c9\u03bd \u03ac\u03c1\u03b8\u03c1\u03c9\u03bd {{offences }} \u03c4\u03c0\u03c5 \u03a0.\u039a.\u03b5\u03b9\u03c2 \u03b2\u03ac\u03c1\u03bf\u03c2 \u03c4\u03bf\u03c5 {{surname}} {{name}} \u03c4\u03bf\u03c5 {{ fathername}} \u03ba\u03b1\u03b9 \u03c4\u03b7\u03c2 {{ mothername}} \u03b3\u03b5\u03bd. {{dateOfBirth  }} \u03c3\u03c4\u03b7 {{ placeOfBirth }} \u03ba\u03b1\u03c4. {{ address }}, \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03c2 \u03c4\u03b7\u03bb\u03b5\u03c6\u03ce\u03bd\u03bf\u03c5 {{ tel }}, \u03b7\u03bb\u03b5\u03ba\u03c4\u03c1\u03bf\u03bd\u03b9\u03ba\u03cc \u03c4\u03b1\u03c7\u03c5\u03b4\u03c1\u03bf\u03bc\u03b5\u03af\u03bf {{ email }}, \u03ba\u03ac\u03c4\u03bf\u03c7\u03bf\u03c2 \u03c4\u03bf\u03c5 \u03c5\u03c0 \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03bd {{ DAT }} \u03c0\u03bf\u03c5 \u03b5\u03ba\u03b4\u03cc\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd  {{ issued }} \u03b1\u03c0\u03cc {{place_issued}} \u0391.\u03a6.\u039c : {{ afm }}, \u0394.\u039f.\u03a5 : {{ doy }}\", \"\u0395\u03b9\u03b4\u03b9\u03ba\u03cc\u03c4\u03b5\u03c1\u03b1 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03b5\u03af\u03c4\u03b1\u03b9 \u03cc\u03c4\u03b9 \u03c3\u03c4\u03b9\u03c2 {{dateOfCrime1}}\u03ba\u03b1\u03b9 \u03c0\u03b5\u03c1\u03af \u03ce\u03c1\u03b1 {{ hourOfCrime1}} \u03c3\u03c4\u03bf {{placeOfCrime1}}  \u03c0\u03c1\u03bf\u03ad\u03b2\u03b7\u03c2 \u03c3\u03c4\u03b7\u03bd \u03c0\u03b1\u03c1\u03ac\u03b2\u03b1\u03c3\u03b7 \u03c4\u03bf\u03c5/\u03c4\u03c9\u03bd \u03ac\u03c1\u03b8\u03c1\u03c9\u03bd {{offences1 }} \u03c4\u03c0\u03c5 \u03a0.\u039a.\u03b5\u03b9\u03c2 \u03b2\u03ac\u03c1\u03bf\u03c2 \u03c4\u03bf\u03c5 {{surnamePerperator}}  {{namePerperator}} \u03c4\u03bf\u03c5 {{fathernamePerperator}} \u03ba\u03b1\u03b9 \u03c4\u03b7\u03c2 {{ mothernamePerperator}} \u03b3\u03b5\u03bd. {{dateOfBirthPerperator }} \u03c3\u03c4\u03b7 {{ placeOfBirthPerperator }}\u03ba\u03b1\u03c4.{{ addressPerperator }},\u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03c2 \u03c4\u03b7\u03bb\u03b5\u03c6\u03ce\u03bd\u03bf\u03c5 {{ telPreperator }}, \u03b7\u03bb\u03b5\u03ba\u03c4\u03c1\u03bf\u03bd\u03b9\u03ba\u03cc  \u03c4\u03b1\u03c7\u03c5\u03b4\u03c1\u03bf\u03bc\u03b5\u03af\u03bf\u03c5{{ emailPreperator }}, \u03ba\u03ac\u03c4\u03bf\u03c7\u03bf\u03c2 \u03c4\u03bf\u03c5 \u03c5\u03c0 \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03bd {{ DATperperator }} \u03c0\u03bf\u03c5 \u03b5\u03ba\u03b4\u03cc\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd {{ issuedPerperator }} \u03b1\u03c0\u03cc {{place_issuedPerperator}}\u0391.\u03a6.\u039c : {{ afmPreperator }}, \u0394.\u039f.\u03a5 : {{ doyPrep }}\"],\n  [\"\u0397 \u03c0\u03b1\u03c1\u03bf\u03cd\u03c3\u03b1 \u03ad\u03ba\u03b8\u03b5\u03c3\u03b7 \u03ac\u03c1\u03c7\u03b9\u03c3\u03b5 \u03bd\u03b1 \u03c3\u03c5\u03bd\u03c4\u03ac\u03c3\u03c3\u03b5\u03c4\u03b1\u03b9 \u03c4\u03b7\u03bd  {{start_hour_apologia}} \u03ce\u03c1\u03b1 \u03ba\u03b1\u03b9 \u03c0\u03b5\u03c1\u03b1\u03b9\u03ce\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd {{end_hour_apologia}} \u03ce\u03c1\u03b1. \u0393\u03b9\u03b1 \u03c0\u03af\u03c3\u03c4\u03c9\u03c3\u03b7 \u03c3\u03c5\u03bd\u03c4\u03ac\u03c7\u03b8\u03b7\u03ba\u03b5 \u03b7 \u03c0\u03b1\u03c1\u03bf\u03cd\u03c3\u03b1 \u03ad\u03ba\u03b8\u03b5\u03c3\u03b7 \u03b7 \u03bf\u03c0\u03bf\u03af\u03b1 \u03b1\u03c6\u03bf\u03cd \u03b1\u03bd\u03b1\u03b3\u03bd\u03ce\u03c3\u03b8\u03b7\u03ba\u03b5 \u03ba\u03b1\u03b9 \u03b2\u03b5\u03b2\u03b1\u03b9\u03ce\u03b8\u03b7\u03ba\u03b5, \u03c5\u03c0\u03bf\u03b3\u03c1\u03ac\u03c6\u03b5\u03c4\u03b1\u03b9 \u03c9\u03c2 \u03b1\u03ba\u03bf\u03bb\u03bf\u03cd\u03b8\u03c9\u03c2:\", \"\u0397 \u03c0\u03b1\u03c1\u03bf\u03cd\u03c3\u03b1 \u03ad\u03ba\u03b8\u03b5\u03c3\u03b7 \u03ac\u03c1\u03c7\u03b9\u03c3\u03b5 \u03bd\u03b1 \u03c3\u03c5\u03bd\u03c4\u03ac\u03c3\u03c3\u03b5\u03c4\u03b1\u03b9 \u03c4\u03b7\u03bd  {{start_hour_apologia1}} \u03ce\u03c1\u03b1 \u03ba\u03b1\u03b9 \u03c0\u03b5\u03c1\u03b1\u03b9\u03ce\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd {{end_hour_apologia1}} \u03ce\u03c1\u03b1. \u0393\u03b9\u03b1 \u03c0\u03af\u03c3\u03c4\u03c9\u03c3\u03b7 \u03c3\u03c5\u03bd\u03c4\u03ac\u03c7\u03b8\u03b7\u03ba\u03b5 \u03b7 \u03c0\u03b1\u03c1\u03bf\u03cd\u03c3\u03b1 \u03ad\u03ba\u03b8\u03b5\u03c3\u03b7 \u03b7 \u03bf\u03c0\u03bf\u03af\u03b1 \u03b1\u03c6\u03bf\u03cd \u03b1\u03bd\u03b1\u03b3\u03bd\u03ce\u03c3\u03b8\u03b7\u03ba\u03b5 \u03ba\u03b1\u03b9 \u03b2\u03b5\u03b2\u03b1\u03b9\u03ce\u03b8\u03b7\u03ba\u03b5, \u03c5\u03c0\u03bf\u03b3\u03c1\u03ac\u03c6\u03b5\u03c4\u03b1\u03b9 \u03c9\u03c2 \u03b1\u03ba\u03bf\u03bb\u03bf\u03cd\u03b8\u03c9\u03c2:\"],\n];\n\nlet applied = 0;\nfor (const para of paragraphs.items) {\n  const text = para.text;\n  for (const [oldText, newText] of replacements) {\n    if (text === oldText) {\n      para.insertText(newText, Word.InsertLocation.replace);\n      applied++;\n      break;\n    }\n  }\n}\nawait context.sync();\nreturn \"applied=\" + applied;\n", "ps1": "# Apply the 5 placeholder-renaming edits described in the diff.\n# Each target paragraph holds exactly one run, so replacing the\n# paragraph Range's Text in place keeps paragraph formatting (pPr) intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Index = 2; OldText = '\u03a3\u03c4\u03b7\u03bd {{ place }} \u03c3\u03ae\u03bc\u03b5\u03c1\u03b1 \u03c4\u03b7\u03bd {{ date_num_apologia}} \u03c4\u03bf\u03c5 \u03bc\u03ae\u03bd\u03b1 {{ month_apologia }} \u03c4\u03bf\u03c5 \u03ad\u03c4\u03bf\u03c5\u03c2 {{year_apologia}} \u03b7\u03bc\u03ad\u03c1\u03b1 \u03b5\u03b2\u03b4\u03bf\u03bc\u03ac\u03b4\u03b1\u03c2 {{ day_apologia}} \u03ba\u03b1\u03b9 \u03ce\u03c1\u03b1 {{start_hour_apologia}} \u03b5\u03bd\u03ce\u03c0\u03b9\u03bf\u03bd \u03b5\u03bc\u03bf\u03cd \u03c4\u03bf\u03c5 {{ first_officer }} \u03c4\u03bf\u03c5 {{ policeStation }} \u0398\u03b5\u03c3\u03c3\u03b1\u03bb\u03bf\u03bd\u03af\u03ba\u03b7\u03c2, \u03c0\u03b1\u03c1\u03b9\u03c3\u03c4\u03b1\u03bc\u03ad\u03bd\u03bf\u03c5  \u03ba\u03b1\u03b9 \u03c4\u03bf\u03c5  {{ sec_officer }} \u03c4\u03b7\u03c2 \u03b9\u03b4\u03af\u03b1\u03c2 \u03c5\u03c0\u03b7\u03c1\u03b5\u03c3\u03af\u03b1\u03c2, \u03c0\u03bf\u03c5 \u03c0\u03c1\u03bf\u03c3\u03bb\u03ae\u03c6\u03b8\u03b7\u03ba\u03b5 \u03c9\u03c2 \u0392'' \u0391\u03bd\u03b1\u03ba\u03c1\u03b9\u03c4\u03b9\u03ba\u03cc\u03c2 \u03a5\u03c0\u03ac\u03bb\u03bb\u03b7\u03bb\u03bf\u03c2, \u03b5\u03be\u03b5\u03c4\u03ac\u03b6\u03b5\u03c4\u03b1\u03b9 \u03bf \u03ba\u03b1\u03c4\u03c9\u03c4\u03ad\u03c1\u03c9 \u03c3\u03b7\u03bc\u03b5\u03b9\u03bf\u03cd\u03bc\u03b5\u03bd\u03bf\u03c2 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03bf\u03cd\u03bc\u03b5\u03bd\u03bf\u03c2 ,'; NewText = '\u03a3\u03c4\u03b7\u03bd {{ place1 }} \u03c3\u03ae\u03bc\u03b5\u03c1\u03b1 \u03c4\u03b7\u03bd {{ date_num_apologia1}} \u03c4\u03bf\u03c5 \u03bc\u03ae\u03bd\u03b1 {{ month_apologia1 }} \u03c4\u03bf\u03c5 \u03ad\u03c4\u03bf\u03c5\u03c2 {{year_apologia1}} \u03b7\u03bc\u03ad\u03c1\u03b1 \u03b5\u03b2\u03b4\u03bf\u03bc\u03ac\u03b4\u03b1\u03c2 {{day_apologia1}} \u03ba\u03b1\u03b9 \u03ce\u03c1\u03b1 {{start_hour_apologia1}} \u03b5\u03bd\u03ce\u03c0\u03b9\u03bf\u03bd \u03b5\u03bc\u03bf\u03cd \u03c4\u03bf\u03c5 {{ first_officer }} \u03c4\u03bf\u03c5 {{ policeStation }} \u0398\u03b5\u03c3\u03c3\u03b1\u03bb\u03bf\u03bd\u03af\u03ba\u03b7\u03c2, \u03c0\u03b1\u03c1\u03b9\u03c3\u03c4\u03b1\u03bc\u03ad\u03bd\u03bf\u03c5  \u03ba\u03b1\u03b9 \u03c4\u03bf\u03c5  {{ sec_officer }} \u03c4\u03b7\u03c2 \u03b9\u03b4\u03af\u03b1\u03c2 \u03c5\u03c0\u03b7\u03c1\u03b5\u03c3\u03af\u03b1\u03c2, \u03c0\u03bf\u03c5 \u03c0\u03c1\u03bf\u03c3\u03bb\u03ae\u03c6\u03b8\u03b7\u03ba\u03b5 \u03c9\u03c2 \u0392'' \u0391\u03bd\u03b1\u03ba\u03c1\u03b9\u03c4\u03b9\u03ba\u03cc\u03c2 \u03a5\u03c0\u03ac\u03bb\u03bb\u03b7\u03bb\u03bf\u03c2, \u03b5\u03be\u03b5\u03c4\u03ac\u03b6\u03b5\u03c4\u03b1\u03b9 \u03bf \u03ba\u03b1\u03c4\u03c9\u03c4\u03ad\u03c1\u03c9 \u03c3\u03b7\u03bc\u03b5\u03b9\u03bf\u03cd\u03bc\u03b5\u03bd\u03bf\u03c2 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03bf\u03cd\u03bc\u03b5\u03bd\u03bf\u03c2 ,' }\n    @{ Index = 4; OldText = '\u0391\u03a0\u039f\u039a\u03a1\u0399\u03a3\u0397 : {{surnamePerperator}}  {{namePerperator}} \u03c4\u03bf\u03c5 {{fathernamePerperator}} \u03ba\u03b1\u03b9 \u03c4\u03b7\u03c2 {{ mothernamePerperator}} \u03b3\u03b5\u03bd. {{dateOfBirthPerperator }} \u03c3\u03c4\u03b7 {{ placeOfBirthPerperator }} \u03ba\u03b1\u03c4.{{ addressPerperator }},\u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03c2 \u03c4\u03b7\u03bb\u03b5\u03c6\u03ce\u03bd\u03bf\u03c5 {{ telPreperator }}, \u03b7\u03bb\u03b5\u03ba\u03c4\u03c1\u03bf\u03bd\u03b9\u03ba\u03cc \u03c4\u03b1\u03c7\u03c5\u03b4\u03c1\u03bf\u03bc\u03b5\u03af\u03bf\u03c5{{ emailPreperator }}, \u03ba\u03ac\u03c4\u03bf\u03c7\u03bf\u03c2 \u03c4\u03bf\u03c5 \u03c5\u03c0 \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03bd {{ DATperperator }} \u03c0\u03bf\u03c5 \u03b5\u03ba\u03b4\u03cc\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd {{ issuedPerperator }} \u03b1\u03c0\u03cc {{place_issuedPerperator}}\u0391.\u03a6.\u039c : {{ afmPreperator }}, \u0394.\u039f.\u03a5 : {{ doyPrep }}, \u03b1\u03c0\u03cc \u03c4\u03bf\u03bd {{officer_arrest }} \u03c5\u03c0\u03b7\u03c1\u03b5\u03c4\u03ce\u03bd \u03c3\u03c4\u03bf  {{  policeStation }}'; NewText = '\u0391\u03a0\u039f\u039a\u03a1\u0399\u03a3\u0397 : {{surname}} {{name}} \u03c4\u03bf\u03c5 {{ fathername}} \u03ba\u03b1\u03b9 \u03c4\u03b7\u03c2 {{ mothername}} \u03b3\u03b5\u03bd. {{dateOfBirth  }} \u03c3\u03c4\u03b7 {{ placeOfBirth }} \u03ba\u03b1\u03c4. {{ address }}, \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03c2 \u03c4\u03b7\u03bb\u03b5\u03c6\u03ce\u03bd\u03bf\u03c5 {{ tel }}, \u03b7\u03bb\u03b5\u03ba\u03c4\u03c1\u03bf\u03bd\u03b9\u03ba\u03cc \u03c4\u03b1\u03c7\u03c5\u03b4\u03c1\u03bf\u03bc\u03b5\u03af\u03bf {{ email }}, \u03ba\u03ac\u03c4\u03bf\u03c7\u03bf\u03c2 \u03c4\u03bf\u03c5 \u03c5\u03c0 \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03bd {{ DAT }} \u03c0\u03bf\u03c5 \u03b5\u03ba\u03b4\u03cc\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd  {{ issued }} \u03b1\u03c0\u03cc {{place_issued}} \u0391.\u03a6.\u039c : {{ afm }}, \u0394.\u039f.\u03a5 : {{ doy }}, \u03b1\u03c0\u03cc \u03c4\u03bf\u03bd {{officer_arrest }} \u03c5\u03c0\u03b7\u03c1\u03b5\u03c4\u03ce\u03bd \u03c3\u03c4\u03bf  {{  policeStation }}' }\n    @{ Index = 5; OldText = '\u0395\u03bd\u03c4\u03b1\u03cd\u03b8\u03b1 \u03b3\u03bd\u03c9\u03c1\u03af\u03c3\u03b1\u03bc\u03b5 \u03c3\u03c4\u03bf\u03bd \u03b5\u03be\u03b5\u03c4\u03b1\u03b6\u03cc\u03bc\u03b5\u03bd\u03bf \u03cc\u03c4\u03b9 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03b5\u03af\u03c4\u03b1\u03b9 \u03b3\u03b9\u03b1 \u03c0\u03b1\u03c1\u03ac\u03b2\u03b1\u03c3\u03b7 \u03c4\u03bf\u03c5 \u03ac\u03c1\u03b8\u03c1\u03bf\u03c5 {{offences }} \u03c4\u03bf\u03c5 \u03a0.\u039a.'; NewText = '\u0395\u03bd\u03c4\u03b1\u03cd\u03b8\u03b1 \u03b3\u03bd\u03c9\u03c1\u03af\u03c3\u03b1\u03bc\u03b5 \u03c3\u03c4\u03bf\u03bd \u03b5\u03be\u03b5\u03c4\u03b1\u03b6\u03cc\u03bc\u03b5\u03bd\u03bf \u03cc\u03c4\u03b9 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03b5\u03af\u03c4\u03b1\u03b9 \u03b3\u03b9\u03b1 \u03c0\u03b1\u03c1\u03ac\u03b2\u03b1\u03c3\u03b7 \u03c4\u03bf\u03c5 \u03ac\u03c1\u03b8\u03c1\u03bf\u03c5 {{offences_text1}} \u03c4\u03bf\u03c5 \u03a0.\u039a.' }\n    @{ Index = 6; OldText = '\u0395\u03b9\u03b4\u03b9\u03ba\u03cc\u03c4\u03b5\u03c1\u03b1 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03b5\u03af\u03c4\u03b1\u03b9 \u03cc\u03c4\u03b9 \u03c3\u03c4\u03b9\u03c2 {{dateOfCrime}}\u03ba\u03b1\u03b9 \u03c0\u03b5\u03c1\u03af \u03ce\u03c1\u03b1 {{ hourOfCrime}} \u03c3\u03c4\u03bf {{placeOfCrime}}  \u03c0\u03c1\u03bf\u03ad\u03b2\u03b7\u03c2 \u03c3\u03c4\u03b7\u03bd \u03c0\u03b1\u03c1\u03ac\u03b2\u03b1\u03c3\u03b7 \u03c4\u03bf\u03c5/\u03c4\u03c9\u03bd \u03ac\u03c1\u03b8\u03c1\u03c9\u03bd {{offences }} \u03c4\u03c0\u03c5 \u03a0.\u039a.\u03b5\u03b9\u03c2 \u03b2\u03ac\u03c1\u03bf\u03c2 \u03c4\u03bf\u03c5 {{surname}} {{name}} \u03c4\u03bf\u03c5 {{ fathername}} \u03ba\u03b1\u03b9 \u03c4\u03b7\u03c2 {{ mothername}} \u03b3\u03b5\u03bd. {{dateOfBirth  }} \u03c3\u03c4\u03b7 {{ placeOfBirth }} \u03ba\u03b1\u03c4. {{ address }}, \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03c2 \u03c4\u03b7\u03bb\u03b5\u03c6\u03ce\u03bd\u03bf\u03c5 {{ tel }}, \u03b7\u03bb\u03b5\u03ba\u03c4\u03c1\u03bf\u03bd\u03b9\u03ba\u03cc \u03c4\u03b1\u03c7\u03c5\u03b4\u03c1\u03bf\u03bc\u03b5\u03af\u03bf {{ email }}, \u03ba\u03ac\u03c4\u03bf\u03c7\u03bf\u03c2 \u03c4\u03bf\u03c5 \u03c5\u03c0 \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03bd {{ DAT }} \u03c0\u03bf\u03c5 \u03b5\u03ba\u03b4\u03cc\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd  {{ issued }} \u03b1\u03c0\u03cc {{place_issued}} \u0391.\u03a6.\u039c : {{ afm }}, \u0394.\u039f.\u03a5 : {{ doy }}'; NewText = '\u0395\u03b9\u03b4\u03b9\u03ba\u03cc\u03c4\u03b5\u03c1\u03b1 \u03ba\u03b1\u03c4\u03b7\u03b3\u03bf\u03c1\u03b5\u03af\u03c4\u03b1\u03b9 \u03cc\u03c4\u03b9 \u03c3\u03c4\u03b9\u03c2 {{dateOfCrime1}}\u03ba\u03b1\u03b9 \u03c0\u03b5\u03c1\u03af \u03ce\u03c1\u03b1 {{ hourOfCrime1}} \u03c3\u03c4\u03bf {{placeOfCrime1}}  \u03c0\u03c1\u03bf\u03ad\u03b2\u03b7\u03c2 \u03c3\u03c4\u03b7\u03bd \u03c0\u03b1\u03c1\u03ac\u03b2\u03b1\u03c3\u03b7 \u03c4\u03bf\u03c5/\u03c4\u03c9\u03bd \u03ac\u03c1\u03b8\u03c1\u03c9\u03bd {{offences1 }} \u03c4\u03c0\u03c5 \u03a0.\u039a.\u03b5\u03b9\u03c2 \u03b2\u03ac\u03c1\u03bf\u03c2 \u03c4\u03bf\u03c5 {{surnamePerperator}}  {{namePerperator}} \u03c4\u03bf\u03c5 {{fathernamePerperator}} \u03ba\u03b1\u03b9 \u03c4\u03b7\u03c2 {{ mothernamePerperator}} \u03b3\u03b5\u03bd. {{dateOfBirthPerperator }} \u03c3\u03c4\u03b7 {{ placeOfBirthPerperator }}\u03ba\u03b1\u03c4.{{ addressPerperator }},\u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03c2 \u03c4\u03b7\u03bb\u03b5\u03c6\u03ce\u03bd\u03bf\u03c5 {{ telPreperator }}, \u03b7\u03bb\u03b5\u03ba\u03c4\u03c1\u03bf\u03bd\u03b9\u03ba\u03cc  \u03c4\u03b1\u03c7\u03c5\u03b4\u03c1\u03bf\u03bc\u03b5\u03af\u03bf\u03c5{{ emailPreperator }}, \u03ba\u03ac\u03c4\u03bf\u03c7\u03bf\u03c2 \u03c4\u03bf\u03c5 \u03c5\u03c0 \u03b1\u03c1\u03b9\u03b8\u03bc\u03cc\u03bd {{ DATperperator }} \u03c0\u03bf\u03c5 \u03b5\u03ba\u03b4\u03cc\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd {{ issuedPerperator }} \u03b1\u03c0\u03cc {{place_issuedPerperator}}\u0391.\u03a6.\u039c : {{ afmPreperator }}, \u0394.\u039f.\u03a5 : {{ doyPrep }}' }\n    @{ Index = 19; OldText = '\u0397 \u03c0\u03b1\u03c1\u03bf\u03cd\u03c3\u03b1 \u03ad\u03ba\u03b8\u03b5\u03c3\u03b7 \u03ac\u03c1\u03c7\u03b9\u03c3\u03b5 \u03bd\u03b1 \u03c3\u03c5\u03bd\u03c4\u03ac\u03c3\u03c3\u03b5\u03c4\u03b1\u03b9 \u03c4\u03b7\u03bd  {{start_hour_apologia}} \u03ce\u03c1\u03b1 \u03ba\u03b1\u03b9 \u03c0\u03b5\u03c1\u03b1\u03b9\u03ce\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd {{end_hour_apologia}} \u03ce\u03c1\u03b1. \u0393\u03b9\u03b1 \u03c0\u03af\u03c3\u03c4\u03c9\u03c3\u03b7 \u03c3\u03c5\u03bd\u03c4\u03ac\u03c7\u03b8\u03b7\u03ba\u03b5 \u03b7 \u03c0\u03b1\u03c1\u03bf\u03cd\u03c3\u03b1 \u03ad\u03ba\u03b8\u03b5\u03c3\u03b7 \u03b7 \u03bf\u03c0\u03bf\u03af\u03b1 \u03b1\u03c6\u03bf\u03cd \u03b1\u03bd\u03b1\u03b3\u03bd\u03ce\u03c3\u03b8\u03b7\u03ba\u03b5 \u03ba\u03b1\u03b9 \u03b2\u03b5\u03b2\u03b1\u03b9\u03ce\u03b8\u03b7\u03ba\u03b5, \u03c5\u03c0\u03bf\u03b3\u03c1\u03ac\u03c6\u03b5\u03c4\u03b1\u03b9 \u03c9\u03c2 \u03b1\u03ba\u03bf\u03bb\u03bf\u03cd\u03b8\u03c9\u03c2:'; NewText = '\u0397 \u03c0\u03b1\u03c1\u03bf\u03cd\u03c3\u03b1 \u03ad\u03ba\u03b8\u03b5\u03c3\u03b7 \u03ac\u03c1\u03c7\u03b9\u03c3\u03b5 \u03bd\u03b1 \u03c3\u03c5\u03bd\u03c4\u03ac\u03c3\u03c3\u03b5\u03c4\u03b1\u03b9 \u03c4\u03b7\u03bd  {{start_hour_apologia1}} \u03ce\u03c1\u03b1 \u03ba\u03b1\u03b9 \u03c0\u03b5\u03c1\u03b1\u03b9\u03ce\u03b8\u03b7\u03ba\u03b5 \u03c4\u03b7\u03bd {{end_hour_apologia1}} \u03ce\u03c1\u03b1. \u0393\u03b9\u03b1 \u03c0\u03af\u03c3\u03c4\u03c9\u03c3\u03b7 \u03c3\u03c5\u03bd\u03c4\u03ac\u03c7\u03b8\u03b7\u03ba\u03b5 \u03b7 \u03c0\u03b1\u03c1\u03bf\u03cd\u03c3\u03b1 \u03ad\u03ba\u03b8\u03b5\u03c3\u03b7 \u03b7 \u03bf\u03c0\u03bf\u03af\u03b1 \u03b1\u03c6\u03bf\u03cd \u03b1\u03bd\u03b1\u03b3\u03bd\u03ce\u03c3\u03b8\u03b7\u03ba\u03b5 \u03ba\u03b1\u03b9 \u03b2\u03b5\u03b2\u03b1\u03b9\u03ce\u03b8\u03b7\u03ba\u03b5, \u03c5\u03c0\u03bf\u03b3\u03c1\u03ac\u03c6\u03b5\u03c4\u03b1\u03b9 \u03c9\u03c2 \u03b1\u03ba\u03bf\u03bb\u03bf\u03cd\u03b8\u03c9\u03c2:' }\n)\n\nforeach ($r in $replacements) {\n    $p = $d.Paragraphs($r.Index)\n    $rng = $p.Range\n    # Paragraph.Range includes the trailing paragraph mark; trim it off\n    # so only the visible text is compared/replaced (mark stays intact).\n    $rngText = $rng.Text\n    if ($rngText.EndsWith([char]13)) {\n        $rngText = $rngText.Substring(0, $rngText.Length - 1)\n    }\n    if ($rngText -ne $r.OldText) {\n        throw \"Paragraph $($r.Index) text did not match expected old text.\"\n    }\n    $rng.MoveEnd(1, -1) | Out-Null\n    $rng.Text = $r.NewText\n}\n\nWrite-Output \"done\"\n"}
